$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 74 - this shifts the existing rows 74:111 down to 75:112,
# preserving their content and formatting (including the date style on column D).
$ws.Rows("74").Insert()

# Populate the newly inserted row 74 with a new weekly price observation for Jengibre.
$ws.Range("A74").Value = 9
$ws.Range("B74").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C74").Value = "Metropolitana"
$ws.Range("D74").Value = 44879
$ws.Range("E74").Value = 13
$ws.Range("F74").Value = 100114007
$ws.Range("G74").Value = "Jengibre"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 380
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = 14526
$ws.Range("N74").Value = "$/caja 13 kilos"
$ws.Range("O74").Value = "Perú"
$ws.Range("P74").Value = 1117
$ws.Range("Q74").Value = 13
$ws.Range("R74").Value = "Hortaliza"

# Make sure the D74 cell keeps the same number format style used by the rest of
# the date column (style index 2 in the original workbook).
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
